$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 480.875
$ws.Range("I28").Value = 405.85715
$ws.Range("J28").Value = 1006
$ws.Range("K28").Value = 405.85715
$ws.Range("L28").Value = 1006
$ws.Range("M28").Value = 79.14285000000001
$ws.Range("N28").Value = -1976
$ws.Range("H31").Value = 4666.6665
$ws.Range("I31").Value = 4666.6665
$ws.Range("K31").Value = 13999.9995
$ws.Range("M31").Value = -13769.9995
$ws.Range("H33").Value = 742.44446
$ws.Range("I33").Value = 363.33334
$ws.Range("K33").Value = 363.33334
$ws.Range("M33").Value = -134.33334
$ws.Range("H38").Value = 7418.5625
$ws.Range("I38").Value = 6726.467
$ws.Range("K38").Value = 20179.401
$ws.Range("M38").Value = -19807.401
$ws.Range("H70").Value = 1484.3334
$ws.Range("I70").Value = 1450
$ws.Range("J70").Value = 1501.5
$ws.Range("K70").Value = 4350
$ws.Range("L70").Value = 4504.5
$ws.Range("M70").Value = -4080
$ws.Range("N70").Value = -5044.5
$ws.Range("H73").Value = 1484.3334
$ws.Range("I73").Value = 1450
$ws.Range("J73").Value = 1501.5
$ws.Range("K73").Value = 4350
$ws.Range("L73").Value = 4504.5
$ws.Range("M73").Value = -3414
$ws.Range("N73").Value = -6376.5
$ws.Range("H112").Value = 3990
$ws.Range("J112").Value = 3990
$ws.Range("L112").Value = 11970
$ws.Range("N112").Value = -14186
$ws.Range("H131").Value = 1472.4
$ws.Range("I131").Value = 1472.4
$ws.Range("K131").Value = 4417.200000000001
$ws.Range("M131").Value = 622.7999999999993
$ws.Range("H135").Value = 1298.0714
$ws.Range("I135").Value = 1657.7
$ws.Range("J135").Value = 399
$ws.Range("K135").Value = 14919.3
$ws.Range("L135").Value = 3591
$ws.Range("M135").Value = -12384.3
$ws.Range("N135").Value = -8661
$ws.Range("H137").Value = 2682.1428
$ws.Range("I137").Value = 2755.1667
$ws.Range("K137").Value = 8265.500100000001
$ws.Range("M137").Value = -5715.500100000001
$ws.Range("H141").Value = 638.8570999999999
$ws.Range("I141").Value = 578.6667
$ws.Range("J141").Value = 1000
$ws.Range("K141").Value = 1736.0001
$ws.Range("L141").Value = 3000
$ws.Range("M141").Value = 3443.9999
$ws.Range("N141").Value = -13360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1268.75
$ws.Range("I2").Value = 1257.1428
$ws.Range("J2").Value = 1350
$ws.Range("K2").Value = 1257.1428
$ws.Range("L2").Value = 1350
$ws.Range("M2").Value = -1144.1428
$ws.Range("N2").Value = -1576
$ws.Range("H37").Value = 10000
$ws.Range("I37").Value = 10000
$ws.Range("K37").Value = 10000
$ws.Range("M37").Value = -9727
$ws.Range("H45").Value = 751.8889
$ws.Range("I45").Value = 408.375
$ws.Range("K45").Value = 408.375
$ws.Range("M45").Value = -31.375
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H102").Value = 1092.5
$ws.Range("I102").Value = 1092.5
$ws.Range("K102").Value = 1092.5
$ws.Range("M102").Value = 529.5
$ws.Range("H116").Value = 1268.75
$ws.Range("I116").Value = 1257.1428
$ws.Range("J116").Value = 1350
$ws.Range("K116").Value = 1257.1428
$ws.Range("L116").Value = 1350
$ws.Range("M116").Value = 1036.8572
$ws.Range("N116").Value = -5938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1268.75
$ws.Range("I3").Value = 1257.1428
$ws.Range("J3").Value = 1350
$ws.Range("K3").Value = 1257.1428
$ws.Range("L3").Value = 1350
$ws.Range("M3").Value = -1143.1428
$ws.Range("N3").Value = -1578
$ws.Range("H105").Value = 1620
$ws.Range("I105").Value = 1025
$ws.Range("K105").Value = 1025
$ws.Range("M105").Value = 722

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1445
$ws.Range("I16").Value = 1426.6666
$ws.Range("K16").Value = 1426.6666
$ws.Range("M16").Value = -1139.6666
$ws.Range("H29").Value = 4990
$ws.Range("I29").Value = 4990
$ws.Range("K29").Value = 4990
$ws.Range("M29").Value = -4697
$ws.Range("H31").Value = 2495.75
$ws.Range("I31").Value = 2406.1
$ws.Range("K31").Value = 2406.1
$ws.Range("M31").Value = -2111.1
$ws.Range("H34").Value = 2495.75
$ws.Range("I34").Value = 2406.1
$ws.Range("K34").Value = 2406.1
$ws.Range("M34").Value = -2204.1
$ws.Range("H74").Value = 61000
$ws.Range("J74").Value = 61000
$ws.Range("L74").Value = 61000
$ws.Range("N74").Value = -62748
$ws.Range("H77").Value = 61000
$ws.Range("J77").Value = 61000
$ws.Range("L77").Value = 183000
$ws.Range("N77").Value = -191736
$ws.Range("H113").Value = 1445
$ws.Range("I113").Value = 1426.6666
$ws.Range("K113").Value = 1426.6666
$ws.Range("M113").Value = 743.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5300
$ws.Range("I4").Value = 5300
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 15900
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -15788
$ws.Range("N4").ClearContents()
$ws.Range("H34").Value = 4433.222
$ws.Range("J34").Value = 4924.875
$ws.Range("L34").Value = 14774.625
$ws.Range("N34").Value = -14942.625
$ws.Range("H39").Value = 5500
$ws.Range("J39").Value = 5291.6665
$ws.Range("L39").Value = 15874.9995
$ws.Range("N39").Value = -16462.9995
$ws.Range("H55").Value = 3864.9
$ws.Range("J55").Value = 4921.2856
$ws.Range("L55").Value = 14763.8568
$ws.Range("N55").Value = -15117.8568
$ws.Range("H75").Value = 808.3333
$ws.Range("J75").Value = 808.3333
$ws.Range("L75").Value = 2424.9999
$ws.Range("N75").Value = -4420.9999
$ws.Range("H78").Value = 808.3333
$ws.Range("J78").Value = 808.3333
$ws.Range("L78").Value = 7274.9997
$ws.Range("N78").Value = -17258.9997
$ws.Range("H114").Value = 218.8
$ws.Range("I114").Value = 346.66666
$ws.Range("J114").Value = 27
$ws.Range("K114").Value = 1039.99998
$ws.Range("L114").Value = 81
$ws.Range("M114").Value = 2214.00002
$ws.Range("N114").Value = -6589
$ws.Range("H117").Value = 4083.3333
$ws.Range("J117").Value = 4083.3333
$ws.Range("L117").Value = 12249.9999
$ws.Range("N117").Value = -19133.9999
$ws.Range("H140").Value = 96.75
$ws.Range("I140").Value = 96.75
$ws.Range("K140").Value = 290.25
$ws.Range("M140").Value = 4889.75
$ws.Range("H141").Value = 2749.75
$ws.Range("J141").Value = 2699.5
$ws.Range("L141").Value = 8098.5
$ws.Range("N141").Value = -18458.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H46").Value = 1314.1428
$ws.Range("I46").Value = 1314.1428
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1314.1428
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1158.1428
$ws.Range("N46").ClearContents()
$ws.Range("H113").Value = 5093.1113
$ws.Range("I113").Value = 5093.1113
$ws.Range("K113").Value = 5093.1113
$ws.Range("M113").Value = -2923.1113

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 2000
$ws.Range("I34").Value = 2000
$ws.Range("K34").Value = 2000
$ws.Range("M34").Value = -1828
$ws.Range("H35").Value = 3337.8
$ws.Range("I35").Value = 1873
$ws.Range("J35").Value = 5535
$ws.Range("K35").Value = 1873
$ws.Range("L35").Value = 5535
$ws.Range("M35").Value = -1537
$ws.Range("N35").Value = -6207
$ws.Range("H93").Value = 3854
$ws.Range("I93").Value = 3758.7144
$ws.Range("J93").Value = 4187.5
$ws.Range("K93").Value = 3758.7144
$ws.Range("L93").Value = 4187.5
$ws.Range("M93").Value = -2510.7144
$ws.Range("N93").Value = -6683.5
$ws.Range("H100").Value = 850
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H127").Value = 42667
$ws.Range("J127").Value = 42667
$ws.Range("L127").Value = 42667
$ws.Range("N127").Value = -52587

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 2011.75
$ws.Range("I122").Value = 1442.1666
$ws.Range("J122").Value = 2581.3333
$ws.Range("K122").Value = 4326.4998
$ws.Range("L122").Value = 7743.999899999999
$ws.Range("M122").Value = -1876.4998
$ws.Range("N122").Value = -12643.9999
